# Auto-generated COM-interop script to apply the 'Generate Report for Handback' edit.
$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item(1)
$wsZhCn = $wb.Worksheets.Item(2)
$wsDeDe = $wb.Worksheets.Item(3)

# --- $wsOverview: rebuild hyperlinks with refreshed file identifiers ---
$wsOverview.Range("A1").Hyperlinks.Delete()
$wsOverview.Hyperlinks.Add($wsOverview.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/940db9130cd23a455ec77d9b56bc0d191305caca/e2e/5fff7693-9672-48a5-a098-040185dcd281.md", "", "", "3804a772-8d6b-4bb4-9cdf-503441cd50cc.md") | Out-Null
$wsOverview.Hyperlinks.Add($wsOverview.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/940db9130cd23a455ec77d9b56bc0d191305caca/e2e/d38901d7-938c-410b-b6cc-a1d01d19b6b0.md", "", "", "ffff03f5a4ce-8f31-4d45-8634-fefbd8e017a5.md") | Out-Null

# --- $wsZhCn: rebuild hyperlinks with refreshed file identifiers ---
$wsZhCn.Range("A1").Hyperlinks.Delete()
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/940db9130cd23a455ec77d9b56bc0d191305caca/e2e/5fff7693-9672-48a5-a098-040185dcd281.md", "", "", "3804a772-8d6b-4bb4-9cdf-503441cd50cc.md") | Out-Null
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("B2"), "https://github.com/OpenLocalizationTest/oltest/blob/940db9130cd23a455ec77d9b56bc0d191305caca/e2e/5fff7693-9672-48a5-a098-040185dcd281.md", "", "", ".md") | Out-Null
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("D2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/7bf4cd131e7678f00e2aef823110b5f50a7dfe95/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/5fff7693-9672-48a5-a098-040185dcd281.49795b3c02d3501ff87a411f67f0c70fca8a892b.zh-cn.xlf", "", "", "3804a772-8d6b-4bb4-9cdf-503441cd50cc.6c8417bee85d1fa480419431ae17850f1352be01.zh-cn.xlf") | Out-Null
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("F2"), "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/a110c469b5d59fed78fcad7a28beaed5cfe70b62/e2e/5fff7693-9672-48a5-a098-040185dcd281.md", "", "", "3804a772-8d6b-4bb4-9cdf-503441cd50cc.md") | Out-Null
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("G2"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/901f1aca164fa7709ae603e3231303adb17e12a8/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/5fff7693-9672-48a5-a098-040185dcd281.49795b3c02d3501ff87a411f67f0c70fca8a892b.zh-cn.xlf", "", "", "3804a772-8d6b-4bb4-9cdf-503441cd50cc.6c8417bee85d1fa480419431ae17850f1352be01.zh-cn.xlf") | Out-Null
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/940db9130cd23a455ec77d9b56bc0d191305caca/e2e/d38901d7-938c-410b-b6cc-a1d01d19b6b0.md", "", "", "ffff03f5a4ce-8f31-4d45-8634-fefbd8e017a5.md") | Out-Null
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("B3"), "https://github.com/OpenLocalizationTest/oltest/blob/940db9130cd23a455ec77d9b56bc0d191305caca/e2e/d38901d7-938c-410b-b6cc-a1d01d19b6b0.md", "", "", ".md") | Out-Null
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("D3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/7bf4cd131e7678f00e2aef823110b5f50a7dfe95/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/d38901d7-938c-410b-b6cc-a1d01d19b6b0.0fa3a3fdfd0601929490cb78620888f992fa0d30.zh-cn.xlf", "", "", "3804a772-8d6b-4bb4-9cdf-503441cd50cc.6c8417bee85d1fa480419431ae17850f1352be01.zh-cn.xlf") | Out-Null
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("F3"), "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/a110c469b5d59fed78fcad7a28beaed5cfe70b62/e2e/d38901d7-938c-410b-b6cc-a1d01d19b6b0.md", "", "", "ffff03f5a4ce-8f31-4d45-8634-fefbd8e017a5.md") | Out-Null
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("G3"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/901f1aca164fa7709ae603e3231303adb17e12a8/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/d38901d7-938c-410b-b6cc-a1d01d19b6b0.0fa3a3fdfd0601929490cb78620888f992fa0d30.zh-cn.xlf", "", "", "3804a772-8d6b-4bb4-9cdf-503441cd50cc.6c8417bee85d1fa480419431ae17850f1352be01.zh-cn.xlf") | Out-Null

# --- $wsDeDe: rebuild hyperlinks with refreshed file identifiers ---
$wsDeDe.Range("A1").Hyperlinks.Delete()
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/940db9130cd23a455ec77d9b56bc0d191305caca/e2e/5fff7693-9672-48a5-a098-040185dcd281.md", "", "", "3804a772-8d6b-4bb4-9cdf-503441cd50cc.md") | Out-Null
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("B2"), "https://github.com/OpenLocalizationTest/oltest/blob/940db9130cd23a455ec77d9b56bc0d191305caca/e2e/5fff7693-9672-48a5-a098-040185dcd281.md", "", "", ".md") | Out-Null
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("D2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/d8592f88ee95f04f1657f9d3b259c77105383fdf/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/5fff7693-9672-48a5-a098-040185dcd281.49795b3c02d3501ff87a411f67f0c70fca8a892b.de-de.xlf", "", "", "3804a772-8d6b-4bb4-9cdf-503441cd50cc.6c8417bee85d1fa480419431ae17850f1352be01.de-de.xlf") | Out-Null
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("F2"), "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/fcecdcbbd0f4ac12c9b647fee4b802b5611c9d10/e2e/5fff7693-9672-48a5-a098-040185dcd281.md", "", "", "3804a772-8d6b-4bb4-9cdf-503441cd50cc.md") | Out-Null
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("G2"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/d916d71ad2dbce3b7ffd9c4c58a3aba5582bc23d/ol-handback/OpenLocalizationTestOrg/oltest.de-de/ci/ht/5fff7693-9672-48a5-a098-040185dcd281.49795b3c02d3501ff87a411f67f0c70fca8a892b.de-de.xlf", "", "", "3804a772-8d6b-4bb4-9cdf-503441cd50cc.6c8417bee85d1fa480419431ae17850f1352be01.de-de.xlf") | Out-Null
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/940db9130cd23a455ec77d9b56bc0d191305caca/e2e/d38901d7-938c-410b-b6cc-a1d01d19b6b0.md", "", "", "ffff03f5a4ce-8f31-4d45-8634-fefbd8e017a5.md") | Out-Null
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("B3"), "https://github.com/OpenLocalizationTest/oltest/blob/940db9130cd23a455ec77d9b56bc0d191305caca/e2e/d38901d7-938c-410b-b6cc-a1d01d19b6b0.md", "", "", ".md") | Out-Null
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("D3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/d8592f88ee95f04f1657f9d3b259c77105383fdf/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/d38901d7-938c-410b-b6cc-a1d01d19b6b0.0fa3a3fdfd0601929490cb78620888f992fa0d30.de-de.xlf", "", "", "3804a772-8d6b-4bb4-9cdf-503441cd50cc.6c8417bee85d1fa480419431ae17850f1352be01.de-de.xlf") | Out-Null
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("F3"), "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/fcecdcbbd0f4ac12c9b647fee4b802b5611c9d10/e2e/d38901d7-938c-410b-b6cc-a1d01d19b6b0.md", "", "", "ffff03f5a4ce-8f31-4d45-8634-fefbd8e017a5.md") | Out-Null
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("G3"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/d916d71ad2dbce3b7ffd9c4c58a3aba5582bc23d/ol-handback/OpenLocalizationTestOrg/oltest.de-de/ci/ht/d38901d7-938c-410b-b6cc-a1d01d19b6b0.0fa3a3fdfd0601929490cb78620888f992fa0d30.de-de.xlf", "", "", "3804a772-8d6b-4bb4-9cdf-503441cd50cc.6c8417bee85d1fa480419431ae17850f1352be01.de-de.xlf") | Out-Null

# --- timestamp columns (plain text cells, not hyperlinked) ---
$wsZhCn.Range("E2").Value = "2016-03-13 09:05:47"
$wsZhCn.Range("H2").Value = "2016-03-13 09:06:04"
$wsZhCn.Range("E3").Value = "2016-03-13 09:05:47"
$wsZhCn.Range("H3").Value = "2016-03-13 09:06:04"

$wsDeDe.Range("E2").Value = "2016-03-13 09:05:50"
$wsDeDe.Range("H2").Value = "2016-03-13 09:06:10"
$wsDeDe.Range("E3").Value = "2016-03-13 09:05:50"
$wsDeDe.Range("H3").Value = "2016-03-13 09:06:10"

